# Latest df -> control of RT (200 ms, 3 SD) and accuracy (0.5 - chance level)
# Negate the FPM (column D) and WPM (column E) values for rows 2-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 41; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value2 = -1 * $dCell.Value2

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value2 = -1 * $eCell.Value2
}
